# Add season record columns (Wins, Losses, Ties) to the TEX_2001 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Copy the header formatting used by the other header cells (e.g. AC1) so
# the new headers match the existing bold/centered/bordered look.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every data row (2 through 53).
$wins = 73
$losses = 89
$ties = 0

for ($row = 2; $row -le 53; $row++) {
    $ws.Cells.Item($row, 30).Value2 = $wins    # column AD
    $ws.Cells.Item($row, 31).Value2 = $losses  # column AE
    $ws.Cells.Item($row, 32).Value2 = $ties    # column AF
}
